{"js": "// Replace each \"old\" equation text with its \"new\" equation text.\n// Every old string is unique in the document and none of the new\n// strings collide with any old string, so a straightforward\n// search-and-replace (old -> new), run once per pair, is unambiguous\n// and safe to do in any order.\nconst replacements = [\n  ['80-24=56', '28+70=98'],\n  ['41+5=46', '39+56=95'],\n  ['61+2=63', '3+0=3'],\n  ['84-38=46', '4+77=81'],\n  ['20+40=60', '18+32=50'],\n  ['41+7=48', '12-9=3'],\n  ['75-21=54', '5+66=71'],\n  ['6+32=38', '89-13=76'],\n  ['73-45=28', '92-65=27'],\n  ['75-50=25', '83+4=87'],\n  ['43-23=20', '26+16=42'],\n  ['16+63=79', '83-12=71'],\n  ['34+35=69', '46+50=96'],\n  ['87-59=28', '42+10=52'],\n  ['76-25=51', '38+22=60'],\n  ['70-35=35', '83+6=89'],\n  ['73-44=29', '38+0=38'],\n  ['70-67=3', '38+40=78'],\n  ['5+85=90', '32+0=32'],\n  ['35+44=79', '68-42=26'],\n  ['88-79=9', '99-55=44'],\n  ['47+46=93', '64+5=69'],\n  ['96-31=65', '71-40=31'],\n  ['69-48=21', '34-26=8'],\n  ['46-10=36', '8-6=2'],\n  ['88-32=56', '44+49=93'],\n  ['88-9=79', '59+15=74'],\n  ['59+25=84', '5+75=80'],\n  ['78-41=37', '39+6=45'],\n  ['47-42=5', '36+58=94'],\n  ['99-53=46', '30-12=18'],\n  ['13+50=63', '55-52=3'],\n  ['97-85=12', '91-55=36'],\n  ['92-42=50', '21+61=82'],\n  ['92-86=6', '15+53=68'],\n  ['72-65=7', '80-54=26'],\n  ['97-80=17', '44-13=31'],\n  ['94-39=55', '6+37=43'],\n  ['53+40=93', '24-22=2'],\n  ['93-53=40', '48-14=34'],\n  ['25+72=97', '72-55=17'],\n  ['10+57=67', '94-43=51'],\n  ['90+2=92', '14+25=39'],\n  ['58-0=58', '89-18=71'],\n  ['94-79=15', '99-27=72'],\n  ['94-25=69', '93-25=68'],\n  ['64-7=57', '26+52=78'],\n  ['21-7=14', '23-12=11'],\n  ['80+18=98', '21+19=40'],\n  ['3+1=4', '46-42=4'],\n  ['83-74=9', '33-19=14'],\n  ['44-12=32', '63-32=31'],\n  ['90+5=95', '18+37=55'],\n  ['88-67=21', '84-34=50'],\n  ['13+16=29', '35+15=50'],\n  ['81-79=2', '98-37=61'],\n  ['44+26=70', '53+28=81'],\n  ['99-71=28', '90-36=54'],\n  ['50+29=79', '31+7=38'],\n  ['61+0=61', '78+7=85'],\n  ['91-83=8', '82-20=62'],\n  ['51+41=92', '75-51=24'],\n  ['35+24=59', '76-65=11'],\n  ['14+8=22', '54+5=59'],\n  ['61+18=79', '27-13=14'],\n  ['7+35=42', '6+77=83'],\n  ['46+12=58', '44+22=66'],\n  ['23-15=8', '74-38=36'],\n  ['23+38=61', '34-22=12'],\n  ['78-59=19', '7+59=66'],\n  ['81+17=98', '65-14=51'],\n  ['75-10=65', '61-35=26'],\n  ['91-5=86', '92-7=85'],\n  ['14+10=24', '46-7=39'],\n  ['79-37=42', '17+66=83'],\n  ['94-22=72', '76-67=9'],\n  ['47+40=87', '42+4=46'],\n  ['17+10=27', '61+25=86'],\n  ['85+8=93', '50+34=84'],\n  ['33+51=84', '93-58=35'],\n  ['49+2=51', '33+30=63'],\n  ['30+46=76', '74-40=34'],\n  ['68+10=78', '27+3=30'],\n  ['29+66=95', '14+0=14'],\n  ['32+10=42', '43-33=10'],\n  ['50-5=45', '82-66=16'],\n  ['61+30=91', '65-18=47'],\n  ['28-1=27', '40+42=82'],\n  ['52-26=26', '6+20=26'],\n  ['36+43=79', '48-23=25'],\n  ['78-68=10', '92-45=47'],\n  ['47-7=40', '11+17=28'],\n  ['10+51=61', '1+95=96'],\n  ['24-10=14', '37+21=58'],\n  ['85-6=79', '95-48=47'],\n  ['28-10=18', '61+15=76'],\n  ['48+14=62', '81+4=85'],\n  ['88-28=60', '86-50=36'],\n  ['75-5=70', '98-30=68'],\n  ['31+28=59', '74+17=91']\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load('items');\n  await context.sync();\n\n  if (found.items.length === 0) {\n    continue; // already applied / not present - skip safely\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"old\" equation text with its \"new\" equation text using\n# Word's Find/Replace (Find.Execute with ReplaceWith = wdReplaceAll).\n# Every old string is unique in the document and none of the new\n# strings collide with any old string, so running these replacements\n# in any order, one old->new substitution per pair, is unambiguous.\n$replacements = @(\n  @(\"80-24=56\", \"28+70=98\"),\n  @(\"41+5=46\", \"39+56=95\"),\n  @(\"61+2=63\", \"3+0=3\"),\n  @(\"84-38=46\", \"4+77=81\"),\n  @(\"20+40=60\", \"18+32=50\"),\n  @(\"41+7=48\", \"12-9=3\"),\n  @(\"75-21=54\", \"5+66=71\"),\n  @(\"6+32=38\", \"89-13=76\"),\n  @(\"73-45=28\", \"92-65=27\"),\n  @(\"75-50=25\", \"83+4=87\"),\n  @(\"43-23=20\", \"26+16=42\"),\n  @(\"16+63=79\", \"83-12=71\"),\n  @(\"34+35=69\", \"46+50=96\"),\n  @(\"87-59=28\", \"42+10=52\"),\n  @(\"76-25=51\", \"38+22=60\"),\n  @(\"70-35=35\", \"83+6=89\"),\n  @(\"73-44=29\", \"38+0=38\"),\n  @(\"70-67=3\", \"38+40=78\"),\n  @(\"5+85=90\", \"32+0=32\"),\n  @(\"35+44=79\", \"68-42=26\"),\n  @(\"88-79=9\", \"99-55=44\"),\n  @(\"47+46=93\", \"64+5=69\"),\n  @(\"96-31=65\", \"71-40=31\"),\n  @(\"69-48=21\", \"34-26=8\"),\n  @(\"46-10=36\", \"8-6=2\"),\n  @(\"88-32=56\", \"44+49=93\"),\n  @(\"88-9=79\", \"59+15=74\"),\n  @(\"59+25=84\", \"5+75=80\"),\n  @(\"78-41=37\", \"39+6=45\"),\n  @(\"47-42=5\", \"36+58=94\"),\n  @(\"99-53=46\", \"30-12=18\"),\n  @(\"13+50=63\", \"55-52=3\"),\n  @(\"97-85=12\", \"91-55=36\"),\n  @(\"92-42=50\", \"21+61=82\"),\n  @(\"92-86=6\", \"15+53=68\"),\n  @(\"72-65=7\", \"80-54=26\"),\n  @(\"97-80=17\", \"44-13=31\"),\n  @(\"94-39=55\", \"6+37=43\"),\n  @(\"53+40=93\", \"24-22=2\"),\n  @(\"93-53=40\", \"48-14=34\"),\n  @(\"25+72=97\", \"72-55=17\"),\n  @(\"10+57=67\", \"94-43=51\"),\n  @(\"90+2=92\", \"14+25=39\"),\n  @(\"58-0=58\", \"89-18=71\"),\n  @(\"94-79=15\", \"99-27=72\"),\n  @(\"94-25=69\", \"93-25=68\"),\n  @(\"64-7=57\", \"26+52=78\"),\n  @(\"21-7=14\", \"23-12=11\"),\n  @(\"80+18=98\", \"21+19=40\"),\n  @(\"3+1=4\", \"46-42=4\"),\n  @(\"83-74=9\", \"33-19=14\"),\n  @(\"44-12=32\", \"63-32=31\"),\n  @(\"90+5=95\", \"18+37=55\"),\n  @(\"88-67=21\", \"84-34=50\"),\n  @(\"13+16=29\", \"35+15=50\"),\n  @(\"81-79=2\", \"98-37=61\"),\n  @(\"44+26=70\", \"53+28=81\"),\n  @(\"99-71=28\", \"90-36=54\"),\n  @(\"50+29=79\", \"31+7=38\"),\n  @(\"61+0=61\", \"78+7=85\"),\n  @(\"91-83=8\", \"82-20=62\"),\n  @(\"51+41=92\", \"75-51=24\"),\n  @(\"35+24=59\", \"76-65=11\"),\n  @(\"14+8=22\", \"54+5=59\"),\n  @(\"61+18=79\", \"27-13=14\"),\n  @(\"7+35=42\", \"6+77=83\"),\n  @(\"46+12=58\", \"44+22=66\"),\n  @(\"23-15=8\", \"74-38=36\"),\n  @(\"23+38=61\", \"34-22=12\"),\n  @(\"78-59=19\", \"7+59=66\"),\n  @(\"81+17=98\", \"65-14=51\"),\n  @(\"75-10=65\", \"61-35=26\"),\n  @(\"91-5=86\", \"92-7=85\"),\n  @(\"14+10=24\", \"46-7=39\"),\n  @(\"79-37=42\", \"17+66=83\"),\n  @(\"94-22=72\", \"76-67=9\"),\n  @(\"47+40=87\", \"42+4=46\"),\n  @(\"17+10=27\", \"61+25=86\"),\n  @(\"85+8=93\", \"50+34=84\"),\n  @(\"33+51=84\", \"93-58=35\"),\n  @(\"49+2=51\", \"33+30=63\"),\n  @(\"30+46=76\", \"74-40=34\"),\n  @(\"68+10=78\", \"27+3=30\"),\n  @(\"29+66=95\", \"14+0=14\"),\n  @(\"32+10=42\", \"43-33=10\"),\n  @(\"50-5=45\", \"82-66=16\"),\n  @(\"61+30=91\", \"65-18=47\"),\n  @(\"28-1=27\", \"40+42=82\"),\n  @(\"52-26=26\", \"6+20=26\"),\n  @(\"36+43=79\", \"48-23=25\"),\n  @(\"78-68=10\", \"92-45=47\"),\n  @(\"47-7=40\", \"11+17=28\"),\n  @(\"10+51=61\", \"1+95=96\"),\n  @(\"24-10=14\", \"37+21=58\"),\n  @(\"85-6=79\", \"95-48=47\"),\n  @(\"28-10=18\", \"61+15=76\"),\n  @(\"48+14=62\", \"81+4=85\"),\n  @(\"88-28=60\", \"86-50=36\"),\n  @(\"75-5=70\", \"98-30=68\"),\n  @(\"31+28=59\", \"74+17=91\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
